# Update gh-pages to output generated at 456a3b4
# Apply updated values to column F ("累计人数"/count) on the "展览" and
# "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 136
$ws1.Range("F6").Value = 95
$ws1.Range("F7").Value = 996
$ws1.Range("F8").Value = 946
$ws1.Range("F10").Value = 102
$ws1.Range("F15").Value = 4229
$ws1.Range("F16").Value = 1241
$ws1.Range("F25").Value = 1510
$ws1.Range("F26").Value = 2453
$ws1.Range("F29").Value = 178
$ws1.Range("F30").Value = 975
$ws1.Range("F33").Value = 21
$ws1.Range("F34").Value = 1406
$ws1.Range("F35").Value = 1986
$ws1.Range("F38").Value = 512

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 136
$ws4.Range("F6").Value = 996
$ws4.Range("F7").Value = 946
$ws4.Range("F10").Value = 102
$ws4.Range("F16").Value = 4229
$ws4.Range("F17").Value = 1241
$ws4.Range("F28").Value = 2453
$ws4.Range("F33").Value = 178
$ws4.Range("F34").Value = 975
$ws4.Range("F37").Value = 1406
$ws4.Range("F38").Value = 1986
$ws4.Range("F43").Value = 512
